# 2025 author list / membership refresh for the STFC collaboration list.
# - Add Marta Sabate-Gilarte as a new collaborator row.
# - Correct Shinji Machida's e-mail address (jiscmail vs collaboration DB).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add Marta Sabate-Gilarte as a new table row (row 18) ------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("B18").Value = "Marta"
$ws.Range("C18").Value = "Sabate-Gilarte"
$ws.Range("F18").Value = "M.Sabate-Gilarte"
$ws.Range("E18").Value = "marta.sabate-gilarte@stfc.ac.uk "
$ws.Range("A18").Value = "Dr."
$ws.Range("D18").Value = "M."
$ws.Range("G18").Value = "STFC-PPD"
$ws.Range("H18").Value = "Particle Physics Department, STFC Rutherford Appleton Laboratory, Harwell Oxford, Didcot, OX11 0QX, UK"
$ws.Range("I18").Value = 0

# --- Fix S. Machida's e-mail address (row 17) ------------------------------
$ws.Range("E17").Value = "shinji.machida@stfc.ac.uk"

# Hyperlink the new e-mail address cell ...
$ws.Hyperlinks.Add($ws.Range("E18"), "mailto:marta.sabate-gilarte@stfc.ac.uk") | Out-Null
$ws.Range("E18").Style = "Hyperlink"

# ... and repoint the existing Machida hyperlink at the corrected address.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$E$17') {
        $hl.Address = "mailto:shinji.machida@stfc.ac.uk"
    }
}

# --- Selection housekeeping (matches the saved view in the workbook) -------
$ws.Range("E17").Select() | Out-Null
